$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 0.5344404915603791
$ws.Range("J4").Value = 0.5197070283324361
$ws.Range("K4").Value = 0.7419874107148443
$ws.Range("L4").Value = 3.074251699119845
